$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-12 Friday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-12-13 Saturday", 2) | Out-Null
$d.Content.Find.Execute("56×31=1736", $true, $true, $false, $false, $false, $true, 1, $false, "48×86=4128", 2) | Out-Null
$d.Content.Find.Execute("70×11=770", $true, $true, $false, $false, $false, $true, 1, $false, "55×48=2640", 2) | Out-Null
$d.Content.Find.Execute("73×32=2336", $true, $true, $false, $false, $false, $true, 1, $false, "68×53=3604", 2) | Out-Null
$d.Content.Find.Execute("89×26=2314", $true, $true, $false, $false, $false, $true, 1, $false, "45×30=1350", 2) | Out-Null
$d.Content.Find.Execute("83×31=2573", $true, $true, $false, $false, $false, $true, 1, $false, "58×19=1102", 2) | Out-Null
$d.Content.Find.Execute("20×11=220", $true, $true, $false, $false, $false, $true, 1, $false, "90×64=5760", 2) | Out-Null
$d.Content.Find.Execute("75×49=3675", $true, $true, $false, $false, $false, $true, 1, $false, "63×89=5607", 2) | Out-Null
$d.Content.Find.Execute("16×13=208", $true, $true, $false, $false, $false, $true, 1, $false, "48×57=2736", 2) | Out-Null
$d.Content.Find.Execute("79×67=5293", $true, $true, $false, $false, $false, $true, 1, $false, "39×87=3393", 2) | Out-Null
$d.Content.Find.Execute("56×72=4032", $true, $true, $false, $false, $false, $true, 1, $false, "17×44=748", 2) | Out-Null
$d.Content.Find.Execute("66×29=1914", $true, $true, $false, $false, $false, $true, 1, $false, "51×12=612", 2) | Out-Null
$d.Content.Find.Execute("16×11=176", $true, $true, $false, $false, $false, $true, 1, $false, "35×68=2380", 2) | Out-Null
$d.Content.Find.Execute("17×62=1054", $true, $true, $false, $false, $false, $true, 1, $false, "54×52=2808", 2) | Out-Null
$d.Content.Find.Execute("88×13=1144", $true, $true, $false, $false, $false, $true, 1, $false, "87×65=5655", 2) | Out-Null
$d.Content.Find.Execute("19×34=646", $true, $true, $false, $false, $false, $true, 1, $false, "75×97=7275", 2) | Out-Null
$d.Content.Find.Execute("26×57=1482", $true, $true, $false, $false, $false, $true, 1, $false, "78×29=2262", 2) | Out-Null
$d.Content.Find.Execute("54×54=2916", $true, $true, $false, $false, $false, $true, 1, $false, "43×67=2881", 2) | Out-Null
$d.Content.Find.Execute("30×16=480", $true, $true, $false, $false, $false, $true, 1, $false, "66×42=2772", 2) | Out-Null
$d.Content.Find.Execute("47×33=1551", $true, $true, $false, $false, $false, $true, 1, $false, "65×90=5850", 2) | Out-Null
$d.Content.Find.Execute("16×96=1536", $true, $true, $false, $false, $false, $true, 1, $false, "25×94=2350", 2) | Out-Null
$d.Content.Find.Execute("39×88=3432", $true, $true, $false, $false, $false, $true, 1, $false, "88×48=4224", 2) | Out-Null
$d.Content.Find.Execute("13×71=923", $true, $true, $false, $false, $false, $true, 1, $false, "93×66=6138", 2) | Out-Null
$d.Content.Find.Execute("56×49=2744", $true, $true, $false, $false, $false, $true, 1, $false, "76×36=2736", 2) | Out-Null
$d.Content.Find.Execute("16×39=624", $true, $true, $false, $false, $false, $true, 1, $false, "18×41=738", 2) | Out-Null
$d.Content.Find.Execute("48×58=2784", $true, $true, $false, $false, $false, $true, 1, $false, "85×22=1870", 2) | Out-Null
